$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny floating point adjustment to the existing A51 timestamp
$ws.Cells.Item(51, 1).Value = 44364.76966868634

# Append new row 52 with the latest scraped data point
$ws.Cells.Item(52, 1).Value = 44365.76761002875
$ws.Cells.Item(52, 1).NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"

$ws.Cells.Item(52, 2).Value = 78443
$ws.Cells.Item(52, 3).Value = 65954
$ws.Cells.Item(52, 4).Value = 3493
$ws.Cells.Item(52, 5).Value = 2122
$ws.Cells.Item(52, 6).Value = 1504
$ws.Cells.Item(52, 7).Value = 20749
$ws.Cells.Item(52, 8).Value = 1480
$ws.Cells.Item(52, 9).Value = 912
$ws.Cells.Item(52, 10).Value = 197
